$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Restore lost agent-results stats: rows that previously had blank
# "Avg Bounces" / "Max Reward" (cols D/E) for Q-Learning, and were
# missing SARSA results entirely (cols F/G), plus corrected values
# for the rows that already had numbers.

# Row 3 (Learning Rate 5, Discount 0.7, Exploration 10) - was fully blank in D:G
$ws.Range("D3").Value = 12.0397
$ws.Range("E3").Value = 51.8
$ws.Range("F3").Value = 10.3833
$ws.Range("G3").Value = 46.98

# Row 4 (5, 0.8, 10) - values unchanged
$ws.Range("D4").Value = 10.865
$ws.Range("E4").Value = 49
$ws.Range("F4").Value = 10.875
$ws.Range("G4").Value = 35

# Row 5 (10, 0.6, 10) - corrected values
$ws.Range("D5").Value = 11.63715
$ws.Range("E5").Value = 51.39
$ws.Range("F5").Value = 7.23685
$ws.Range("G5").Value = 32.71

# Row 6 (10, 0.7, 10) - corrected values
$ws.Range("D6").Value = 12.0967
$ws.Range("E6").Value = 54.95
$ws.Range("F6").Value = 9.832
$ws.Range("G6").Value = 45.44

# Row 7 (10, 0.8, 10) - was fully blank in D:G
$ws.Range("D7").Value = 12.8845
$ws.Range("E7").Value = 57.52
$ws.Range("F7").Value = 8.6528
$ws.Range("G7").Value = 38.26

# Row 8 (10, 0.8, 20) - corrected values
$ws.Range("D8").Value = 13.0958
$ws.Range("E8").Value = 55.27
$ws.Range("F8").Value = 10.51345
$ws.Range("G8").Value = 47.43

# Row 9 (10, 0.9, 10) - values unchanged
$ws.Range("D9").Value = 13.77
$ws.Range("E9").Value = 57
$ws.Range("F9").Value = 9.6
$ws.Range("G9").Value = 40

# Row 10 (10, 1, 10) - values unchanged
$ws.Range("D10").Value = 8.96
$ws.Range("E10").Value = 45
$ws.Range("F10").Value = 8.055
$ws.Range("G10").Value = 47

# Row 11 (20, 0.8, 10) - corrected values
$ws.Range("D11").Value = 13.152
$ws.Range("E11").Value = 58.1
$ws.Range("F11").Value = 7.90955
$ws.Range("G11").Value = 32.84

# Row 12 (20, 0.9, 10) - was fully blank in D:G
$ws.Range("D12").Value = 14.7683
$ws.Range("E12").Value = 65.58
$ws.Range("F12").Value = 7.7418
$ws.Range("G12").Value = 32.89

# Reflect the author's final selection/scroll position in the sheet view
$ws.Activate() | Out-Null
$ws.Range("G5").Select() | Out-Null
